# Apply updated "想去人数" (interest count) figures across the four sheets
# of the 广州-漫展信息 workbook, per the commit "Update gh-pages to output
# generated at 456a3b4".
#
# Column F on every sheet holds the "想去人数" counter; this edit only
# bumps those counts (plus one carried-over date-serial value on the
# COMICUP row, which also lives in column F). No other cells change.

$wb = $excel.ActiveWorkbook

function Set-FValues($SheetName, $Updates) {
    $ws = $wb.Worksheets.Item($SheetName)
    foreach ($row in $Updates.Keys) {
        $ws.Range("F$row").Value = $Updates[$row]
    }
}

# Sheet "展览" (rId1 / sheet1)
Set-FValues "展览" @{
    2  = 38030
    4  = 643
    5  = 805
    6  = 492
    7  = 381
    9  = 870
    11 = 766
    12 = 599
    13 = 90
    15 = 46
    16 = 698
    17 = 195
    18 = 497
    20 = 1198
    22 = 892
    23 = 2615
    24 = 1098
    25 = 590
    26 = 130
    27 = 1190
    29 = 855
    30 = 79
    31 = 1198
}

# Sheet "演出" (rId2 / sheet2)
Set-FValues "演出" @{
    3  = 460
    10 = 18
}

# Sheet "本地生活" (rId3 / sheet3)
Set-FValues "本地生活" @{
    2 = 680
}

# Sheet "全部类型" (rId4 / sheet4) -- aggregates the rows above
Set-FValues "全部类型" @{
    2  = 680
    3  = 38030
    5  = 643
    6  = 805
    7  = 492
    9  = 381
    11 = 460
    12 = 460
    16 = 870
    18 = 766
    19 = 599
    20 = 90
    25 = 18
    26 = 46
    28 = 698
    29 = 195
    30 = 497
    32 = 1198
    34 = 892
    35 = 2615
    36 = 1098
    37 = 590
    38 = 130
    39 = 1190
    42 = 855
    43 = 79
    44 = 1198
}
